$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = "AXION OLIVOS - Mariano Pelliza 3892 - "
